$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1729489
$ws.Range("I28").Value = 2375166.2
$ws.Range("J28").Value = 7683
$ws.Range("K28").Value = 2375166.2
$ws.Range("L28").Value = 7683
$ws.Range("M28").Value = -2374681.2
$ws.Range("N28").Value = -8653
$ws.Range("H38").Value = 528.4
$ws.Range("I38").Value = 528.4
$ws.Range("K38").Value = 1585.2
$ws.Range("M38").Value = -1213.2
$ws.Range("H41").Value = 1637
$ws.Range("I41").Value = 2209
$ws.Range("J41").Value = 1228.4286
$ws.Range("K41").Value = 2209
$ws.Range("L41").Value = 1228.4286
$ws.Range("M41").Value = -1769
$ws.Range("N41").Value = -2108.4286
$ws.Range("H64").Value = 3328.9333
$ws.Range("I64").Value = 2733.4
$ws.Range("K64").Value = 2733.4
$ws.Range("M64").Value = -2485.4
$ws.Range("H67").Value = 3328.9333
$ws.Range("I67").Value = 2733.4
$ws.Range("K67").Value = 2733.4
$ws.Range("M67").Value = -1875.4
$ws.Range("H127").Value = 2636.0625
$ws.Range("I127").Value = 2409
$ws.Range("J127").Value = 2928
$ws.Range("K127").Value = 7227
$ws.Range("L127").Value = 8784
$ws.Range("M127").Value = -2267
$ws.Range("N127").Value = -18704
$ws.Range("H131").Value = 1536.3448
$ws.Range("I131").Value = 766.0833
$ws.Range("J131").Value = 2080.0588
$ws.Range("K131").Value = 2298.2499
$ws.Range("L131").Value = 6240.176399999999
$ws.Range("M131").Value = 2741.7501
$ws.Range("N131").Value = -16320.1764
$ws.Range("H132").Value = 1431.3549
$ws.Range("I132").Value = 1323.2069
$ws.Range("K132").Value = 3969.620699999999
$ws.Range("M132").Value = -1439.620699999999
$ws.Range("H135").Value = 58824384
$ws.Range("I135").Value = 979.9091
$ws.Range("K135").Value = 8819.1819
$ws.Range("M135").Value = -6284.1819
$ws.Range("H141").Value = 2549357
$ws.Range("I141").Value = 3501433.8
$ws.Range("J141").Value = 10485.333
$ws.Range("K141").Value = 10504301.4
$ws.Range("L141").Value = 31455.999
$ws.Range("M141").Value = -10499121.4
$ws.Range("N141").Value = -41815.999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 5242.4287
$ws.Range("H66").Value = 5242.4287
$ws.Range("H122").Value = 1880
$ws.Range("I122").Value = 1522.4762
$ws.Range("J122").Value = 3131.3333
$ws.Range("K122").Value = 4567.4286
$ws.Range("L122").Value = 9393.999899999999
$ws.Range("M122").Value = -2117.4286
$ws.Range("N122").Value = -14293.9999
$ws.Range("H132").Value = 1352.6765
$ws.Range("I132").Value = 1093.2069
$ws.Range("J132").Value = 2857.6
$ws.Range("K132").Value = 3279.620699999999
$ws.Range("L132").Value = 8572.799999999999
$ws.Range("M132").Value = -749.6206999999995
$ws.Range("N132").Value = -13632.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3916.6667
$ws.Range("I20").Value = 4166.6665
$ws.Range("J20").Value = 3666.6667
$ws.Range("K20").Value = 4166.6665
$ws.Range("L20").Value = 3666.6667
$ws.Range("M20").Value = -3919.6665
$ws.Range("N20").Value = -4160.6667
$ws.Range("H134").Value = 5305.7144
$ws.Range("I134").Value = 5478.5186
$ws.Range("J134").Value = 640
$ws.Range("K134").Value = 16435.5558
$ws.Range("L134").Value = 1920
$ws.Range("M134").Value = -13900.5558
$ws.Range("N134").Value = -6990

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3752.8
$ws.Range("I122").Value = 2878
$ws.Range("J122").Value = 5502.4
$ws.Range("K122").Value = 8634
$ws.Range("L122").Value = 16507.2
$ws.Range("M122").Value = -6184
$ws.Range("N122").Value = -21407.2
$ws.Range("H132").Value = 2804.8572
$ws.Range("I132").Value = 2066.476
$ws.Range("K132").Value = 6199.428
$ws.Range("M132").Value = -3669.428
$ws.Range("H134").Value = 1873.3846
$ws.Range("I134").Value = 1705.091
$ws.Range("K134").Value = 5115.272999999999
$ws.Range("M134").Value = -2580.272999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 686.75
$ws.Range("I117").Value = 792
$ws.Range("J117").Value = 651.6667
$ws.Range("K117").Value = 2376
$ws.Range("L117").Value = 1955.0001
$ws.Range("M117").Value = 1066
$ws.Range("N117").Value = -8839.000099999999
$ws.Range("H122").Value = 1962.6666
$ws.Range("I122").Value = 999
$ws.Range("K122").Value = 8991
$ws.Range("M122").Value = -6541

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1317.6
$ws.Range("I113").Value = 1038.5
$ws.Range("J113").Value = 1503.6666
$ws.Range("K113").Value = 1038.5
$ws.Range("L113").Value = 1503.6666
$ws.Range("M113").Value = 1131.5
$ws.Range("N113").Value = -5843.6666
$ws.Range("H122").Value = 1410.2
$ws.Range("I122").Value = 1410.2
$ws.Range("K122").Value = 4230.6
$ws.Range("M122").Value = -1780.6
$ws.Range("H126").Value = 26327.047
$ws.Range("J126").Value = 39402.85
$ws.Range("L126").Value = 118208.55
$ws.Range("N126").Value = -123148.55
$ws.Range("H132").Value = 2901.6843
$ws.Range("I132").Value = 2243.6
$ws.Range("K132").Value = 6730.799999999999
$ws.Range("M132").Value = -4200.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 11456.833
$ws.Range("I16").Value = 11456.833
$ws.Range("K16").Value = 11456.833
$ws.Range("M16").Value = -11286.833
$ws.Range("H22").Value = 2025
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 2025
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 2025
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -2615
$ws.Range("H27").Value = 2025
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 2025
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 2025
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -2239
$ws.Range("H55").Value = 348.375
$ws.Range("I55").Value = 357.5
$ws.Range("J55").Value = 333.16666
$ws.Range("K55").Value = 357.5
$ws.Range("L55").Value = 333.16666
$ws.Range("M55").Value = -184.5
$ws.Range("N55").Value = -679.16666
$ws.Range("H61").Value = 2361.7273
$ws.Range("I61").Value = 1622.375
$ws.Range("J61").Value = 4333.3335
$ws.Range("K61").Value = 1622.375
$ws.Range("L61").Value = 4333.3335
$ws.Range("M61").Value = -1420.375
$ws.Range("N61").Value = -4737.3335
$ws.Range("H82").Value = 1750.625
$ws.Range("J82").Value = 2599.6667
$ws.Range("L82").Value = 2599.6667
$ws.Range("N82").Value = -3321.6667
$ws.Range("H85").Value = 1750.625
$ws.Range("J85").Value = 2599.6667
$ws.Range("L85").Value = 2599.6667
$ws.Range("N85").Value = -5095.6667
$ws.Range("H113").Value = 2361.7273
$ws.Range("I113").Value = 1622.375
$ws.Range("J113").Value = 4333.3335
$ws.Range("K113").Value = 1622.375
$ws.Range("L113").Value = 4333.3335
$ws.Range("M113").Value = 547.625
$ws.Range("N113").Value = -8673.333500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 47599.668
$ws.Range("J123").Value = 47599.668
$ws.Range("L123").Value = 47599.668
$ws.Range("N123").Value = -57399.668
$ws.Range("H132").Value = 2066.926
$ws.Range("I132").Value = 1181.25
$ws.Range("J132").Value = 3355.182
$ws.Range("K132").Value = 3543.75
$ws.Range("L132").Value = 10065.546
$ws.Range("M132").Value = -1013.75
$ws.Range("N132").Value = -15125.546
$ws.Range("H136").Value = 1377.3636
$ws.Range("I136").Value = 1354.0769
$ws.Range("K136").Value = 4062.2307
$ws.Range("M136").Value = -1512.2307
